# Apply the "break out stock.yaml completed" change to the "day" sheet:
#  1. Fix D64:D66 (bsecode) so they are stored as real numbers instead of text.
#  2. Append 10 new data rows (67-76) for the 02/07/2024 11:34:40 batch.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- 1. Convert existing bsecode cells (D64:D66) from text to numeric ---
$ws.Range("D64").Value = 532343
$ws.Range("D65").Value = 533273
$ws.Range("D66").Value = 532321

# --- 2. Append new rows 67-76 ---
# Columns: A=sr, B=nsecode, C=name, D=bsecode(text), E=per_chg, F=close, G=volume, H=timeframe, I=Date Time

$newRows = @(
    @{ r = 67; A = 1;  B = "BAJAJ-AUTO"; C = "Bajaj Auto Limited";                       D = "532977"; E = -1.38; F = 9401.25;             G = 287887;   I = "02/07/2024 11:34:40" },
    @{ r = 68; A = 2;  B = "ABB";        C = "Abb India Limited";                        D = "500002"; E = -0.85; F = 8514.200000000001;  G = 240096;   I = "02/07/2024 11:34:40" },
    @{ r = 69; A = 3;  B = "BAJFINANCE"; C = "Bajaj Finance Limited";                     D = "500034"; E = -1.53; F = 7165.6;              G = 1281027;  I = "02/07/2024 11:34:40" },
    @{ r = 70; A = 4;  B = "HEROMOTOCO"; C = "Hero Motocorp Limited";                     D = "500182"; E = -0.64; F = 5567.1;              G = 689851;   I = "02/07/2024 11:34:40" },
    @{ r = 71; A = 5;  B = "EICHERMOT";  C = "Eicher Motors Limited";                     D = "505200"; E = -0.21; F = 4625.75;             G = 865358;   I = "02/07/2024 11:34:40" },
    @{ r = 72; A = 6;  B = "TITAN";      C = "Titan Company Limited";                     D = "500114"; E = -0.93; F = 3399.65;             G = 878961;   I = "02/07/2024 11:34:40" },
    @{ r = 73; A = 7;  B = "M&M";        C = "Mahindra & Mahindra Limited";               D = "500520"; E = -0.37; F = 2865.15;             G = 1507891;  I = "02/07/2024 11:34:40" },
    @{ r = 74; A = 8;  B = "DLF";        C = "Dlf Limited";                               D = "532868"; E = 0.61;  F = 830.4;               G = 3485222;  I = "02/07/2024 11:34:40" },
    @{ r = 75; A = 9;  B = "ABFRL";      C = "Aditya Birla Fashion And Retail Limited";   D = "535755"; E = 2.19;  F = 329.55;              G = 9171569;  I = "02/07/2024 11:34:40" },
    @{ r = 76; A = 10; B = "TATASTEEL";  C = "Tata Steel Limited";                        D = "500470"; E = 0.27;  F = 174.54;              G = 33359495; I = "02/07/2024 11:34:40" }
)

foreach ($row in $newRows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    # Force bsecode to stay a text value (matches source data) by using the
    # leading-apostrophe convention, same as typing '532977 into Excel.
    $ws.Cells.Item($r, 4).Value = "'" + $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = "day"
    $ws.Cells.Item($r, 9).Value = $row.I
}
